# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E21) listed arrears periods 2407..2412.
# The oldest periods are removed and replaced with newer ones, shifting
# the period values shown for each worker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2411"
$ws.Range("E17").Value = "2410"
$ws.Range("E18").Value = "2409"
$ws.Range("E19").Value = "2408"
$ws.Range("E20").Value = "2407"
